# Scheduled runner: refresh cached Universalis market-price snapshots
# for the Carbuncle data-center profit tables (currentAveragePrice*,
# LevePrice*, LeveProfit* columns) across all job worksheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H57").Value = 10636.3
$ws.Range("J57").Value = 11030.333
$ws.Range("L57").Value = 33090.999
$ws.Range("N57").Value = -34088.999
$ws.Range("H74").Value = 5090.909
$ws.Range("J74").Value = 5250
$ws.Range("L74").Value = 5250
$ws.Range("N74").Value = -7122
$ws.Range("H76").Value = 3607.2896
$ws.Range("I76").Value = 3002.3704
$ws.Range("K76").Value = 3002.3704
$ws.Range("M76").Value = -2687.3704
$ws.Range("H77").Value = 5090.909
$ws.Range("J77").Value = 5250
$ws.Range("L77").Value = 26250
$ws.Range("N77").Value = -35610
$ws.Range("H79").Value = 3607.2896
$ws.Range("I79").Value = 3002.3704
$ws.Range("K79").Value = 3002.3704
$ws.Range("M79").Value = -1910.3704
$ws.Range("H86").Value = 150502.25
$ws.Range("I86").Value = 67336.336
$ws.Range("J86").Value = 400000
$ws.Range("K86").Value = 67336.336
$ws.Range("L86").Value = 400000
$ws.Range("M86").Value = -66213.336
$ws.Range("N86").Value = -402246
$ws.Range("H89").Value = 150502.25
$ws.Range("I89").Value = 67336.336
$ws.Range("J89").Value = 400000
$ws.Range("K89").Value = 336681.68
$ws.Range("L89").Value = 2000000
$ws.Range("M89").Value = -331065.68
$ws.Range("N89").Value = -2011232

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2025
$ws.Range("I2").Value = 2033.3334
$ws.Range("J2").Value = 2000
$ws.Range("K2").Value = 2033.3334
$ws.Range("L2").Value = 2000
$ws.Range("M2").Value = -1920.3334
$ws.Range("N2").Value = -2226
$ws.Range("H45").Value = 1294.9231
$ws.Range("I45").Value = 1166.7273
$ws.Range("K45").Value = 1166.7273
$ws.Range("M45").Value = -789.7273
$ws.Range("H116").Value = 2025
$ws.Range("I116").Value = 2033.3334
$ws.Range("J116").Value = 2000
$ws.Range("K116").Value = 2033.3334
$ws.Range("L116").Value = 2000
$ws.Range("M116").Value = 260.6666
$ws.Range("N116").Value = -6588
$ws.Range("H122").Value = 3444.4443
$ws.Range("I122").Value = 3000
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 9000
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -6550
$ws.Range("N122").Value = -19900

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2025
$ws.Range("I3").Value = 2033.3334
$ws.Range("J3").Value = 2000
$ws.Range("K3").Value = 2033.3334
$ws.Range("L3").Value = 2000
$ws.Range("M3").Value = -1919.3334
$ws.Range("N3").Value = -2228
$ws.Range("H51").Value = 34275
$ws.Range("J51").Value = 34275
$ws.Range("L51").Value = 34275
$ws.Range("N51").Value = -35257
$ws.Range("H82").Value = 6665.2856
$ws.Range("I82").Value = 1976.1666
$ws.Range("J82").Value = 34800
$ws.Range("K82").Value = 1976.1666
$ws.Range("L82").Value = 34800
$ws.Range("M82").Value = -1593.1666
$ws.Range("N82").Value = -35566
$ws.Range("H85").Value = 6665.2856
$ws.Range("I85").Value = 1976.1666
$ws.Range("J85").Value = 34800
$ws.Range("K85").Value = 1976.1666
$ws.Range("L85").Value = 34800
$ws.Range("M85").Value = -650.1666
$ws.Range("N85").Value = -37452
$ws.Range("H86").Value = 2516.48
$ws.Range("I86").Value = 2173.625
$ws.Range("J86").Value = 3126
$ws.Range("K86").Value = 2173.625
$ws.Range("L86").Value = 3126
$ws.Range("M86").Value = -1050.625
$ws.Range("N86").Value = -5372
$ws.Range("H89").Value = 2516.48
$ws.Range("I89").Value = 2173.625
$ws.Range("J89").Value = 3126
$ws.Range("K89").Value = 10868.125
$ws.Range("L89").Value = 15630
$ws.Range("M89").Value = -5252.125
$ws.Range("N89").Value = -26862
$ws.Range("H94").Value = 657.1111
$ws.Range("I94").Value = 778.6667
$ws.Range("J94").Value = 535.55554
$ws.Range("K94").Value = 778.6667
$ws.Range("L94").Value = 535.55554
$ws.Range("M94").Value = -327.6667
$ws.Range("N94").Value = -1437.55554
$ws.Range("H134").Value = 1302.3448
$ws.Range("I134").Value = 1060.381
$ws.Range("J134").Value = 1937.5
$ws.Range("K134").Value = 3181.143
$ws.Range("L134").Value = 5812.5
$ws.Range("M134").Value = -646.143
$ws.Range("N134").Value = -10882.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1002595.4
$ws.Range("I132").Value = 1668419.1
$ws.Range("J132").Value = 3859.8
$ws.Range("K132").Value = 5005257.300000001
$ws.Range("L132").Value = 11579.4
$ws.Range("M132").Value = -5002727.300000001
$ws.Range("N132").Value = -16639.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2926.25
$ws.Range("I80").Value = 2734.1667
$ws.Range("J80").Value = 3502.5
$ws.Range("K80").Value = 2734.1667
$ws.Range("L80").Value = 3502.5
$ws.Range("M80").Value = -1736.1667
$ws.Range("N80").Value = -5498.5
$ws.Range("H83").Value = 2926.25
$ws.Range("I83").Value = 2734.1667
$ws.Range("J83").Value = 3502.5
$ws.Range("K83").Value = 13670.8335
$ws.Range("L83").Value = 17512.5
$ws.Range("M83").Value = -8678.833500000001
$ws.Range("N83").Value = -27496.5
$ws.Range("H122").Value = 103346.08
$ws.Range("I122").Value = 159099.38
$ws.Range("J122").Value = 4229.1113
$ws.Range("K122").Value = 477298.14
$ws.Range("L122").Value = 12687.3339
$ws.Range("M122").Value = -474848.14
$ws.Range("N122").Value = -17587.3339

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H45").Value = 21500
$ws.Range("I45").Value = 10000
$ws.Range("K45").Value = 10000
$ws.Range("M45").Value = -9593
$ws.Range("H53").Value = 8000
$ws.Range("J53").Value = 8000
$ws.Range("L53").Value = 8000
$ws.Range("N53").Value = -9036
$ws.Range("H100").Value = 2910
$ws.Range("I100").Value = 2120
$ws.Range("K100").Value = 2120
$ws.Range("M100").Value = -1579
$ws.Range("H122").Value = 6946684
$ws.Range("I122").Value = 10102579
$ws.Range("K122").Value = 30307737
$ws.Range("M122").Value = -30305287
$ws.Range("H123").Value = 38500
$ws.Range("J123").Value = 38500
$ws.Range("L123").Value = 38500
$ws.Range("N123").Value = -48300

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 42247.6
$ws.Range("J46").Value = 42247.6
$ws.Range("L46").Value = 42247.6
$ws.Range("N46").Value = -42709.6
$ws.Range("H132").Value = 1805.4242
$ws.Range("I132").Value = 1288.24
$ws.Range("J132").Value = 3421.625
$ws.Range("K132").Value = 3864.72
$ws.Range("L132").Value = 10264.875
$ws.Range("M132").Value = -1334.72
$ws.Range("N132").Value = -15324.875
$ws.Range("H134").Value = 42247.6
$ws.Range("J134").Value = 42247.6
$ws.Range("L134").Value = 126742.8
$ws.Range("N134").Value = -131812.8
